# Update label positions and plot size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update a handful of label-position nudge values (x_nudge/y_nudge) ---
$ws.Range("C38").Value = 0.25
$ws.Range("C40").Value = -0.15
$ws.Range("B41").Value = -0.4
$ws.Range("C48").Value = 0.05
$ws.Range("B69").Value = 0.4

# --- Update the plot/view size: scroll the frozen pane down so row 48 is the
#     first visible row beneath the frozen header/column, and move the active
#     selection down to B70. The freeze boundary itself (row 1 / column A)
#     stays the same - only the scrolled viewport and selection change. ---
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B70").Select()
